# "this added last report 30-12-24"
# Updates the daily cash denomination counts on Sheet1 (RSO 02 / RSO 03
# blocks in rows 3-9, and the rows 18-25 reconciliation block).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- RSO 02 (K) / RSO 03 (P) quantity columns, rows 3-9 ---
$ws.Range("K3").Value = 12
$ws.Range("P3").Value = 6

$ws.Range("K4").Value = 23
$ws.Range("P4").Value = 27

$ws.Range("K5").Value = 7
$ws.Range("P5").Value = 8

$ws.Range("K6").Value = 14
$ws.Range("P6").Value = 26

$ws.Range("K7").Value = 4

# K8 had a value (3) and is now blank; P8 is newly populated.
$ws.Range("K8").ClearContents()
$ws.Range("P8").Value = 24

# K9 had a value (1) and is now blank.
$ws.Range("K9").ClearContents()

# --- Rows 18-25 reconciliation block (H / M / R quantity columns) ---
$ws.Range("H18").Value = 1
$ws.Range("M18").Value = 6
$ws.Range("R18").Value = 1

$ws.Range("H19").Value = 5
$ws.Range("M19").Value = 20
$ws.Range("R19").Value = 10

# H20 had a value (47) and is now blank.
$ws.Range("H20").ClearContents()
$ws.Range("M20").Value = 5
$ws.Range("R20").Value = 8

$ws.Range("H21").Value = 4
$ws.Range("M21").Value = 36
$ws.Range("R21").Value = 72

$ws.Range("H22").Value = 3
$ws.Range("R22").Value = 61

$ws.Range("H23").Value = 18
$ws.Range("M23").Value = 9
$ws.Range("R23").Value = 28

$ws.Range("M24").Value = 1
$ws.Range("R24").Value = 34

$ws.Range("H25").Value = 3
$ws.Range("M25").Value = 1

# Restore the selected cell as last left by the author.
$ws.Range("H22").Select()
